# Generate Report for Handoff
#
# The handoff/handback status of the "9f27646a-ddc7-4349-93bb-ac079887f00e.md"
# file moved on from "Handed back: in sync with en-US" to "Ready for handoff",
# with refreshed timestamps and a new error detail noting the handback file
# is stale, for both the zh-cn and de-de locales (and mirrored on the
# Overview roll-up sheet).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fba72dab6752beac315a2860d35a2cbd9d94dbf/e2e/9f27646a-ddc7-4349-93bb-ac079887f00e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9680d151e09b7e48e8f0cfd1d0ec9202c40e96/e2e/9f27646a-ddc7-4349-93bb-ac079887f00e.md."

# --- Overview sheet: row 3 is the 9f27646a-ddc7-4349-93bb-ac079887f00e.md entry ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 10:49:40"

# --- zh-cn sheet: row 3 is the 9f27646a-ddc7-4349-93bb-ac079887f00e.md entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-24 10:49:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the 9f27646a-ddc7-4349-93bb-ac079887f00e.md entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-24 10:49:40"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
